{"js": "const body = context.document.body;\nconst pairs = [\n  [\"2023-08-15 Tuesday\", \"2023-08-16 Wednesday\"],\n  [\"49\u00d777=3773\", \"11\u00d786=946\"],\n  [\"96\u00d754=5184\", \"50\u00d714=700\"],\n  [\"89\u00d717=1513\", \"63\u00d723=1449\"],\n  [\"77\u00d768=5236\", \"59\u00d717=1003\"],\n  [\"91\u00d748=4368\", \"71\u00d718=1278\"],\n  [\"77\u00d721=1617\", \"73\u00d781=5913\"],\n  [\"46\u00d761=2806\", \"19\u00d746=874\"],\n  [\"37\u00d781=2997\", \"68\u00d730=2040\"],\n  [\"94\u00d737=3478\", \"94\u00d731=2914\"],\n  [\"30\u00d720=600\", \"39\u00d755=2145\"],\n  [\"32\u00d719=608\", \"78\u00d750=3900\"],\n  [\"46\u00d788=4048\", \"76\u00d755=4180\"],\n  [\"32\u00d744=1408\", \"24\u00d743=1032\"],\n  [\"88\u00d774=6512\", \"22\u00d714=308\"],\n  [\"77\u00d783=6391\", \"59\u00d732=1888\"],\n  [\"55\u00d761=3355\", \"56\u00d725=1400\"],\n  [\"43\u00d714=602\", \"60\u00d712=720\"],\n  [\"90\u00d774=6660\", \"88\u00d717=1496\"],\n  [\"84\u00d768=5712\", \"20\u00d763=1260\"],\n  [\"36\u00d763=2268\", \"11\u00d793=1023\"],\n  [\"19\u00d792=1748\", \"40\u00d775=3000\"],\n  [\"65\u00d785=5525\", \"24\u00d754=1296\"],\n  [\"89\u00d733=2937\", \"42\u00d784=3528\"],\n  [\"20\u00d776=1520\", \"89\u00d712=1068\"],\n  [\"45\u00d794=4230\", \"64\u00d752=3328\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply text replacements for the two-digit multiplication worksheet update.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2023-08-15 Tuesday\", \"2023-08-16 Wednesday\"),\n  @(\"49\u00d777=3773\", \"11\u00d786=946\"),\n  @(\"96\u00d754=5184\", \"50\u00d714=700\"),\n  @(\"89\u00d717=1513\", \"63\u00d723=1449\"),\n  @(\"77\u00d768=5236\", \"59\u00d717=1003\"),\n  @(\"91\u00d748=4368\", \"71\u00d718=1278\"),\n  @(\"77\u00d721=1617\", \"73\u00d781=5913\"),\n  @(\"46\u00d761=2806\", \"19\u00d746=874\"),\n  @(\"37\u00d781=2997\", \"68\u00d730=2040\"),\n  @(\"94\u00d737=3478\", \"94\u00d731=2914\"),\n  @(\"30\u00d720=600\", \"39\u00d755=2145\"),\n  @(\"32\u00d719=608\", \"78\u00d750=3900\"),\n  @(\"46\u00d788=4048\", \"76\u00d755=4180\"),\n  @(\"32\u00d744=1408\", \"24\u00d743=1032\"),\n  @(\"88\u00d774=6512\", \"22\u00d714=308\"),\n  @(\"77\u00d783=6391\", \"59\u00d732=1888\"),\n  @(\"55\u00d761=3355\", \"56\u00d725=1400\"),\n  @(\"43\u00d714=602\", \"60\u00d712=720\"),\n  @(\"90\u00d774=6660\", \"88\u00d717=1496\"),\n  @(\"84\u00d768=5712\", \"20\u00d763=1260\"),\n  @(\"36\u00d763=2268\", \"11\u00d793=1023\"),\n  @(\"19\u00d792=1748\", \"40\u00d775=3000\"),\n  @(\"65\u00d785=5525\", \"24\u00d754=1296\"),\n  @(\"89\u00d733=2937\", \"42\u00d784=3528\"),\n  @(\"20\u00d776=1520\", \"89\u00d712=1068\"),\n  @(\"45\u00d794=4230\", \"64\u00d752=3328\"),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
